# Auto-generated edit script applying numeric updates to Sheets per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 138.16667
$ws.Range("I2").Value = 138.16667
$ws.Range("K2").Value = 138.16667
$ws.Range("M2").Value = -25.16667000000001
$ws.Range("H9").Value = 1587484.2
$ws.Range("I9").Value = 1587484.2
$ws.Range("K9").Value = 1587484.2
$ws.Range("M9").Value = -1587315.2
$ws.Range("H12").Value = 409.4
$ws.Range("I12").Value = 261.25
$ws.Range("K12").Value = 261.25
$ws.Range("M12").Value = -91.25
$ws.Range("H17").Value = 4145503.8
$ws.Range("J17").Value = 5440242.5
$ws.Range("L17").Value = 16320727.5
$ws.Range("N17").Value = -16321063.5
$ws.Range("H18").Value = 247.14285
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H28").Value = 3149
$ws.Range("I28").Value = 3149
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3149
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -2664
$ws.Range("N28").ClearContents()
$ws.Range("H40").Value = 2426.6
$ws.Range("I40").Value = 2400
$ws.Range("J40").Value = 2499.75
$ws.Range("K40").Value = 2400
$ws.Range("L40").Value = 2499.75
$ws.Range("M40").Value = -2225
$ws.Range("N40").Value = -2849.75
$ws.Range("H58").Value = 1068.25
$ws.Range("I58").Value = 78
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 234
$ws.Range("L58").Value = 24000
$ws.Range("M58").Value = -84
$ws.Range("N58").Value = -24300
$ws.Range("H64").Value = 3688.7778
$ws.Range("I64").Value = 3566.5
$ws.Range("J64").Value = 3933.3333
$ws.Range("K64").Value = 3566.5
$ws.Range("L64").Value = 3933.3333
$ws.Range("M64").Value = -3318.5
$ws.Range("N64").Value = -4429.3333
$ws.Range("H67").Value = 3688.7778
$ws.Range("I67").Value = 3566.5
$ws.Range("J67").Value = 3933.3333
$ws.Range("K67").Value = 3566.5
$ws.Range("L67").Value = 3933.3333
$ws.Range("M67").Value = -2708.5
$ws.Range("N67").Value = -5649.3333
$ws.Range("H74").Value = 8459.200000000001
$ws.Range("I74").Value = 8459.200000000001
$ws.Range("K74").Value = 8459.200000000001
$ws.Range("M74").Value = -7523.200000000001
$ws.Range("H77").Value = 8459.200000000001
$ws.Range("I77").Value = 8459.200000000001
$ws.Range("K77").Value = 42296
$ws.Range("M77").Value = -37616
$ws.Range("H80").Value = 477260.75
$ws.Range("I80").Value = 592.4666999999999
$ws.Range("K80").Value = 1777.4001
$ws.Range("M80").Value = -779.4000999999998
$ws.Range("H83").Value = 477260.75
$ws.Range("I83").Value = 592.4666999999999
$ws.Range("K83").Value = 5332.2003
$ws.Range("M83").Value = -340.2002999999995
$ws.Range("H86").Value = 5183.8335
$ws.Range("I86").Value = 4433.6665
$ws.Range("J86").Value = 5934
$ws.Range("K86").Value = 4433.6665
$ws.Range("L86").Value = 5934
$ws.Range("M86").Value = -3310.6665
$ws.Range("N86").Value = -8180
$ws.Range("H89").Value = 5183.8335
$ws.Range("I89").Value = 4433.6665
$ws.Range("J89").Value = 5934
$ws.Range("K89").Value = 22168.3325
$ws.Range("L89").Value = 29670
$ws.Range("M89").Value = -16552.3325
$ws.Range("N89").Value = -40902
$ws.Range("H98").Value = 1286.8667
$ws.Range("I98").Value = 1366.9166
$ws.Range("K98").Value = 1366.9166
$ws.Range("M98").Value = 131.0834
$ws.Range("H100").Value = 36111.066
$ws.Range("I100").Value = 41185.96
$ws.Range("J100").Value = 3124.25
$ws.Range("K100").Value = 41185.96
$ws.Range("L100").Value = 3124.25
$ws.Range("M100").Value = -40644.96
$ws.Range("N100").Value = -4206.25
$ws.Range("H101").Value = 637.1111
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H103").Value = 1052
$ws.Range("I103").Value = 986
$ws.Range("K103").Value = 2958
$ws.Range("M103").Value = -2372
$ws.Range("H111").Value = 1280.3334
$ws.Range("I111").Value = 1280.3334
$ws.Range("K111").Value = 3841.0002
$ws.Range("M111").Value = -774.0001999999999
$ws.Range("H116").Value = 9666.643
$ws.Range("I116").Value = 7314
$ws.Range("J116").Value = 12019.286
$ws.Range("K116").Value = 7314
$ws.Range("L116").Value = 12019.286
$ws.Range("M116").Value = -3872
$ws.Range("N116").Value = -18903.286
$ws.Range("H122").Value = 1286.8667
$ws.Range("I122").Value = 1366.9166
$ws.Range("K122").Value = 4100.7498
$ws.Range("M122").Value = -1650.7498
$ws.Range("H129").Value = 2215.8333
$ws.Range("I129").Value = 708.9
$ws.Range("J129").Value = 4099.5
$ws.Range("K129").Value = 2126.7
$ws.Range("L129").Value = 12298.5
$ws.Range("M129").Value = 2873.3
$ws.Range("N129").Value = -22298.5
$ws.Range("H132").Value = 2127.8235
$ws.Range("I132").Value = 2099.4546
$ws.Range("K132").Value = 6298.3638
$ws.Range("M132").Value = -3768.3638
$ws.Range("H135").Value = 2286.2
$ws.Range("I135").Value = 2438.7856
$ws.Range("J135").Value = 150
$ws.Range("K135").Value = 21949.0704
$ws.Range("L135").Value = 1350
$ws.Range("M135").Value = -19414.0704
$ws.Range("N135").Value = -6420
$ws.Range("H137").Value = 8758.290999999999
$ws.Range("I137").Value = 3570.8438
$ws.Range("K137").Value = 10712.5314
$ws.Range("M137").Value = -8162.5314
$ws.Range("H138").Value = 2897
$ws.Range("I138").Value = 2842.4
$ws.Range("J138").Value = 2929.76
$ws.Range("K138").Value = 8527.200000000001
$ws.Range("L138").Value = 8789.280000000001
$ws.Range("M138").Value = -3387.200000000001
$ws.Range("N138").Value = -19069.28
$ws.Range("H141").Value = 2627.5293
$ws.Range("I141").Value = 2564.6
$ws.Range("J141").Value = 3099.5
$ws.Range("K141").Value = 7693.799999999999
$ws.Range("L141").Value = 9298.5
$ws.Range("M141").Value = -2513.799999999999
$ws.Range("N141").Value = -19658.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2836.5908
$ws.Range("I32").Value = 1278.6487
$ws.Range("J32").Value = 11071.429
$ws.Range("K32").Value = 1278.6487
$ws.Range("L32").Value = 11071.429
$ws.Range("M32").Value = -991.6487
$ws.Range("N32").Value = -11645.429
$ws.Range("H61").Value = 4388.3794
$ws.Range("I61").Value = 2570.7273
$ws.Range("J61").Value = 10101
$ws.Range("K61").Value = 2570.7273
$ws.Range("L61").Value = 10101
$ws.Range("M61").Value = -2358.7273
$ws.Range("N61").Value = -10525
$ws.Range("H74").Value = 2224.516
$ws.Range("I74").Value = 1758.55
$ws.Range("J74").Value = 3071.7273
$ws.Range("K74").Value = 1758.55
$ws.Range("L74").Value = 3071.7273
$ws.Range("M74").Value = -884.55
$ws.Range("N74").Value = -4819.7273
$ws.Range("H77").Value = 2224.516
$ws.Range("I77").Value = 1758.55
$ws.Range("J77").Value = 3071.7273
$ws.Range("K77").Value = 8792.75
$ws.Range("L77").Value = 15358.6365
$ws.Range("M77").Value = -4424.75
$ws.Range("N77").Value = -24094.6365
$ws.Range("H97").Value = 796.9655
$ws.Range("I97").Value = 807.6070999999999
$ws.Range("J97").Value = 499
$ws.Range("K97").Value = 807.6070999999999
$ws.Range("L97").Value = 499
$ws.Range("M97").Value = -311.6070999999999
$ws.Range("N97").Value = -1491
$ws.Range("H102").Value = 5383.2383
$ws.Range("I102").Value = 5452.4
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 5452.4
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -3830.4
$ws.Range("N102").Value = -7244
$ws.Range("H110").Value = 2019.7307
$ws.Range("I110").Value = 2003.0476
$ws.Range("J110").Value = 2089.8
$ws.Range("K110").Value = 2003.0476
$ws.Range("L110").Value = 2089.8
$ws.Range("M110").Value = 41.9523999999999
$ws.Range("N110").Value = -6179.8
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678
$ws.Range("H122").Value = 3125.111
$ws.Range("I122").Value = 1946.5
$ws.Range("J122").Value = 5482.3335
$ws.Range("K122").Value = 5839.5
$ws.Range("L122").Value = 16447.0005
$ws.Range("M122").Value = -3389.5
$ws.Range("N122").Value = -21347.0005
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 6979.5405
$ws.Range("I132").Value = 7117.8613
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 21353.5839
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -18823.5839
$ws.Range("N132").Value = -11060
$ws.Range("H136").Value = 4388.3794
$ws.Range("I136").Value = 2570.7273
$ws.Range("J136").Value = 10101
$ws.Range("K136").Value = 7712.1819
$ws.Range("L136").Value = 30303
$ws.Range("M136").Value = -5162.1819
$ws.Range("N136").Value = -35403

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 47799.2
$ws.Range("J92").Value = 47799.2
$ws.Range("L92").Value = 47799.2
$ws.Range("N92").Value = -52791.2
$ws.Range("H94").Value = 355.80646
$ws.Range("I94").Value = 334.31033
$ws.Range("K94").Value = 334.31033
$ws.Range("M94").Value = 116.68967
$ws.Range("H99").Value = 4317.1113
$ws.Range("I99").Value = 4494.25
$ws.Range("K99").Value = 4494.25
$ws.Range("M99").Value = -2996.25
$ws.Range("H107").Value = 2995
$ws.Range("I107").Value = 2995
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2995
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1075
$ws.Range("N107").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 303.27274
$ws.Range("I7").Value = 248.55556
$ws.Range("K7").Value = 248.55556
$ws.Range("M7").Value = -135.55556
$ws.Range("H31").Value = 3038.7437
$ws.Range("J31").Value = 4969.2856
$ws.Range("L31").Value = 4969.2856
$ws.Range("N31").Value = -5559.2856
$ws.Range("H34").Value = 3038.7437
$ws.Range("J34").Value = 4969.2856
$ws.Range("L34").Value = 4969.2856
$ws.Range("N34").Value = -5373.2856
$ws.Range("H58").Value = 2989.4285
$ws.Range("I58").Value = 1563.0476
$ws.Range("J58").Value = 7268.5713
$ws.Range("K58").Value = 1563.0476
$ws.Range("L58").Value = 7268.5713
$ws.Range("M58").Value = -1360.0476
$ws.Range("N58").Value = -7674.5713
$ws.Range("H86").Value = 4999.857
$ws.Range("I86").Value = 4999.8335
$ws.Range("K86").Value = 4999.8335
$ws.Range("M86").Value = -3876.8335
$ws.Range("H89").Value = 4999.857
$ws.Range("I89").Value = 4999.8335
$ws.Range("K89").Value = 24999.1675
$ws.Range("M89").Value = -19383.1675
$ws.Range("H95").Value = 33974.625
$ws.Range("J95").Value = 33974.625
$ws.Range("L95").Value = 33974.625
$ws.Range("N95").Value = -39466.625
$ws.Range("H105").Value = 1176.909
$ws.Range("I105").Value = 1266.875
$ws.Range("K105").Value = 1266.875
$ws.Range("M105").Value = 480.125
$ws.Range("H134").Value = 9339.767
$ws.Range("I134").Value = 7954.9546
$ws.Range("J134").Value = 13148
$ws.Range("K134").Value = 23864.8638
$ws.Range("L134").Value = 39444
$ws.Range("M134").Value = -21329.8638
$ws.Range("N134").Value = -44514
$ws.Range("H136").Value = 2989.4285
$ws.Range("I136").Value = 1563.0476
$ws.Range("J136").Value = 7268.5713
$ws.Range("K136").Value = 4689.142800000001
$ws.Range("L136").Value = 21805.7139
$ws.Range("M136").Value = -2139.142800000001
$ws.Range("N136").Value = -26905.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 292.66666
$ws.Range("I2").Value = 29.25
$ws.Range("J2").Value = 819.5
$ws.Range("K2").Value = 175.5
$ws.Range("L2").Value = 4917
$ws.Range("M2").Value = -62.5
$ws.Range("N2").Value = -5143
$ws.Range("H5").Value = 2123.7354
$ws.Range("J5").Value = 2625.652
$ws.Range("L5").Value = 7876.956
$ws.Range("N5").Value = -8100.956
$ws.Range("H7").Value = 280.9
$ws.Range("I7").Value = 395.6
$ws.Range("J7").Value = 166.2
$ws.Range("K7").Value = 1186.8
$ws.Range("L7").Value = 498.6
$ws.Range("M7").Value = -1074.8
$ws.Range("N7").Value = -722.5999999999999
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("H104").Value = 7876.7
$ws.Range("I104").Value = 4935.4
$ws.Range("J104").Value = 10818
$ws.Range("K104").Value = 14806.2
$ws.Range("L104").Value = 32454
$ws.Range("M104").Value = -12185.2
$ws.Range("N104").Value = -37696
$ws.Range("H117").Value = 979.5
$ws.Range("I117").Value = 1071.75
$ws.Range("J117").Value = 795
$ws.Range("K117").Value = 3215.25
$ws.Range("L117").Value = 2385
$ws.Range("M117").Value = 226.75
$ws.Range("N117").Value = -9269
$ws.Range("H122").Value = 10001559
$ws.Range("I122").Value = 2631
$ws.Range("J122").Value = 14286814
$ws.Range("K122").Value = 23679
$ws.Range("L122").Value = 128581326
$ws.Range("M122").Value = -21229
$ws.Range("N122").Value = -128586226
$ws.Range("H129").Value = 2191.5386
$ws.Range("J129").Value = 2427.75
$ws.Range("L129").Value = 7283.25
$ws.Range("N129").Value = -17283.25
$ws.Range("H131").Value = 20683.326
$ws.Range("I131").Value = 167206
$ws.Range("J131").Value = 1571.674
$ws.Range("K131").Value = 501618
$ws.Range("L131").Value = 4715.022
$ws.Range("M131").Value = -496578
$ws.Range("N131").Value = -14795.022
$ws.Range("H134").Value = 1488.4
$ws.Range("I134").Value = 1488.4
$ws.Range("K134").Value = 4465.200000000001
$ws.Range("M134").Value = 604.7999999999993
$ws.Range("H135").Value = 2123.7354
$ws.Range("J135").Value = 2625.652
$ws.Range("L135").Value = 23630.868
$ws.Range("N135").Value = -28700.868
$ws.Range("H136").Value = 1118.2
$ws.Range("I136").Value = 899.75
$ws.Range("K136").Value = 2699.25
$ws.Range("M136").Value = 2400.75
$ws.Range("H137").Value = 3038.9092
$ws.Range("J137").Value = 3072
$ws.Range("L137").Value = 9216
$ws.Range("N137").Value = -19416
$ws.Range("H138").Value = 2488
$ws.Range("I138").Value = 2488
$ws.Range("K138").Value = 7464
$ws.Range("M138").Value = -2324

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 36672
$ws.Range("I33").Value = 35017
$ws.Range("J33").Value = 37499.5
$ws.Range("K33").Value = 35017
$ws.Range("L33").Value = 37499.5
$ws.Range("M33").Value = -34765
$ws.Range("N33").Value = -38003.5
$ws.Range("H46").Value = 54995
$ws.Range("J46").Value = 54995
$ws.Range("L46").Value = 54995
$ws.Range("N46").Value = -55307
$ws.Range("H57").Value = 29979.076
$ws.Range("I57").Value = 9959
$ws.Range("J57").Value = 38876.89
$ws.Range("K57").Value = 9959
$ws.Range("L57").Value = 38876.89
$ws.Range("M57").Value = -9139
$ws.Range("N57").Value = -40516.89
$ws.Range("H80").Value = 2565.8572
$ws.Range("I80").Value = 2659.6667
$ws.Range("J80").Value = 2495.5
$ws.Range("K80").Value = 2659.6667
$ws.Range("L80").Value = 2495.5
$ws.Range("M80").Value = -1661.6667
$ws.Range("N80").Value = -4491.5
$ws.Range("H83").Value = 2565.8572
$ws.Range("I83").Value = 2659.6667
$ws.Range("J83").Value = 2495.5
$ws.Range("K83").Value = 13298.3335
$ws.Range("L83").Value = 12477.5
$ws.Range("M83").Value = -8306.333500000001
$ws.Range("N83").Value = -22461.5
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 79015.266
$ws.Range("J113").Value = 2813.9092
$ws.Range("L113").Value = 2813.9092
$ws.Range("N113").Value = -7153.9092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3474.1428
$ws.Range("I7").Value = 1831.25
$ws.Range("J7").Value = 5664.6665
$ws.Range("K7").Value = 1831.25
$ws.Range("L7").Value = 5664.6665
$ws.Range("M7").Value = -1719.25
$ws.Range("N7").Value = -5888.6665
$ws.Range("H16").Value = 2083.9285
$ws.Range("I16").Value = 1474.8462
$ws.Range("J16").Value = 10002
$ws.Range("K16").Value = 1474.8462
$ws.Range("L16").Value = 10002
$ws.Range("M16").Value = -1304.8462
$ws.Range("N16").Value = -10342
$ws.Range("H46").Value = 1591.6923
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1591.6923
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1591.6923
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1967.6923
$ws.Range("H63").Value = 32570.666
$ws.Range("J63").Value = 27085
$ws.Range("L63").Value = 27085
$ws.Range("N63").Value = -28583
$ws.Range("H66").Value = 32570.666
$ws.Range("J66").Value = 27085
$ws.Range("L66").Value = 81255
$ws.Range("N66").Value = -88743
$ws.Range("H93").Value = 4364.778
$ws.Range("I93").Value = 4441.75
$ws.Range("J93").Value = 3749
$ws.Range("K93").Value = 4441.75
$ws.Range("L93").Value = 3749
$ws.Range("M93").Value = -3193.75
$ws.Range("N93").Value = -6245
$ws.Range("H97").Value = 38687.125
$ws.Range("J97").Value = 38687.125
$ws.Range("L97").Value = 38687.125
$ws.Range("N97").Value = -40669.125
$ws.Range("H100").Value = 1014881.94
$ws.Range("I100").Value = 6170.5713
$ws.Range("J100").Value = 2780126.8
$ws.Range("K100").Value = 6170.5713
$ws.Range("L100").Value = 2780126.8
$ws.Range("M100").Value = -5629.5713
$ws.Range("N100").Value = -2781208.8
$ws.Range("H126").Value = 3474.1428
$ws.Range("I126").Value = 1831.25
$ws.Range("J126").Value = 5664.6665
$ws.Range("K126").Value = 5493.75
$ws.Range("L126").Value = 16993.9995
$ws.Range("M126").Value = -3023.75
$ws.Range("N126").Value = -21933.9995
$ws.Range("H130").Value = 26002.5
$ws.Range("I130").Value = 26002.5
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 26002.5
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -20982.5
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 2851419.2
$ws.Range("I132").Value = 3474426.8
$ws.Range("J132").Value = 3384.1428
$ws.Range("K132").Value = 10423280.4
$ws.Range("L132").Value = 10152.4284
$ws.Range("M132").Value = -10420750.4
$ws.Range("N132").Value = -15212.4284
$ws.Range("H136").Value = 4632104.5
$ws.Range("I136").Value = 6538107
$ws.Range("K136").Value = 19614321
$ws.Range("M136").Value = -19611771

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 74999
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H54").Value = 35266.5
$ws.Range("J54").Value = 35266.5
$ws.Range("L54").Value = 35266.5
$ws.Range("N54").Value = -36306.5
$ws.Range("H81").Value = 1293.7273
$ws.Range("I81").Value = 1293.7273
$ws.Range("K81").Value = 2587.4546
$ws.Range("M81").Value = -1526.4546
$ws.Range("H84").Value = 1293.7273
$ws.Range("I84").Value = 1293.7273
$ws.Range("K84").Value = 12937.273
$ws.Range("M84").Value = -7633.273000000001
$ws.Range("H94").Value = 54994.5
$ws.Range("J94").Value = 54994.5
$ws.Range("L94").Value = 54994.5
$ws.Range("N94").Value = -56796.5
$ws.Range("H122").Value = 5229.9414
$ws.Range("I122").Value = 3664.4614
$ws.Range("J122").Value = 10317.75
$ws.Range("K122").Value = 10993.3842
$ws.Range("L122").Value = 30953.25
$ws.Range("M122").Value = -8543.3842
$ws.Range("N122").Value = -35853.25
$ws.Range("H126").Value = 6352.7144
$ws.Range("I126").Value = 7326.1816
$ws.Range("J126").Value = 2783.3333
$ws.Range("K126").Value = 21978.5448
$ws.Range("L126").Value = 8349.999899999999
$ws.Range("M126").Value = -19508.5448
$ws.Range("N126").Value = -13289.9999
$ws.Range("H127").Value = 44999
$ws.Range("J127").Value = 44999
$ws.Range("L127").Value = 44999
$ws.Range("N127").Value = -54919
$ws.Range("H132").Value = 15377.046
$ws.Range("I132").Value = 10392.383
$ws.Range("J132").Value = 27707.525
$ws.Range("K132").Value = 31177.149
$ws.Range("L132").Value = 83122.57500000001
$ws.Range("M132").Value = -28647.149
$ws.Range("N132").Value = -88182.57500000001
$ws.Range("H136").Value = 825.6842
$ws.Range("I136").Value = 593.7778
$ws.Range("K136").Value = 1781.3334
$ws.Range("M136").Value = 768.6666
